$d = $word.ActiveDocument

function Get-ParagraphByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

function Get-ParagraphContainingText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($text)) {
            return $p
        }
    }
    return $null
}

# ------------------------------------------------------------------
# 1. Delete the trailing "This PDF version is provided under the same
#    license." paragraph completely - including its paragraph mark.
#    (Delete this one first, since it is physically below the other
#    paragraph we are about to delete, keeping earlier positions
#    valid.)
# ------------------------------------------------------------------
$pPdf = Get-ParagraphByText $d "This PDF version is provided under the same license."
if ($pPdf -ne $null) {
    $pPdf.Range.Delete()
}

# ------------------------------------------------------------------
# 2. Delete the "License Information" Heading2 paragraph completely -
#    including its paragraph mark.
# ------------------------------------------------------------------
$pHeading = Get-ParagraphByText $d "License Information"
if ($pHeading -ne $null) {
    $pHeading.Range.Delete()
}

# ------------------------------------------------------------------
# 3. Rewrite the license-description paragraph with the new content,
#    keeping the paragraph's own paragraph mark/run intact.
# ------------------------------------------------------------------
$pLicense = Get-ParagraphContainingText $d "is based on"
$r = $pLicense.Range
$r.End = $r.End - 1        # exclude the paragraph mark
$r.Text = ""                 # clear existing runs, leaving one empty run

$pos = $r.Start

# -- insert all the new text pieces first (no character formatting
#    applied yet, so inserts don't inherit stray formatting) --
$boldText = "unfoldingWord® Translation Questions"
$rIns = $d.Range($pos, $pos)
$rIns.InsertAfter($boldText)
$boldStart = $pos
$boldEnd = $pos + $boldText.Length
$pos = $boldEnd

$plain1 = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. "
$rIns = $d.Range($pos, $pos)
$rIns.InsertAfter($plain1)
$pos = $pos + $plain1.Length

$plain2 = "unfoldingWord® Translation Questions"
$rIns = $d.Range($pos, $pos)
$rIns.InsertAfter($plain2)
$pos = $pos + $plain2.Length

$plain3 = " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from "
$rIns = $d.Range($pos, $pos)
$rIns.InsertAfter($plain3)
$pos = $pos + $plain3.Length

$plain4 = "unfoldingWord® Translation Questions"
$rIns = $d.Range($pos, $pos)
$rIns.InsertAfter($plain4)
$pos = $pos + $plain4.Length

$plain5 = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual"
$rIns = $d.Range($pos, $pos)
$rIns.InsertAfter($plain5)
$pos = $pos + $plain5.Length

# -- as the final step, bold only the first piece ("unfoldingWord®
#    Translation Questions") --
$rBold = $d.Range($boldStart, $boldEnd)
$rBold.Font.Bold = $true
